$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.497.04'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.105.76'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +4.64%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '329.75'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5246'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.14%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4392'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +3.09%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '50.42'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +15.86%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08884'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.38%  '
$ws.Range('E11').Value = '  +2.68%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '24.83'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.56%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.108.27'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +4.64%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.744'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.757'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.85%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '96.57'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.51%  '
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001130'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.53%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06641'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.52%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '19.31'
$ws.Range('D20').Style = "Normal"
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.316'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.78%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '30.560.65'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('E24').Value = '  +3.84%  '
$ws.Range('E25').Value = '  +4.11%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.354.65'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +4.54%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.47'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.629'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +7.25%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '162.25'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '132.83'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.13%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.222'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +7.43%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1072'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.64%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.687'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +23.43%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.235'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.44%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.898'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.78%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '10.27'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +11.42%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02588'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.97%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06753'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.23%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.517'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.87%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '12.73'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.40%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.2279'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.76%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.6922'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.82%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.276'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.88%  '
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '14.09'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.05%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6422'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.91%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.229'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.10%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.634'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.225'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +10.76%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '82.76'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.80%  '
